$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Table 0")

$ws.Range("C2").Value = 0.12180000000000001
$ws.Range("C3").Value = 3.8208000000000002
$ws.Range("C4").Value = 2.911
$ws.Range("C5").Value = 0.49130000000000001
$ws.Range("C6").Value = 3.0341999999999998
$ws.Range("C7").Value = 2.6877
$ws.Range("C8").Value = 2.8475000000000001
$ws.Range("C9").Value = 4.5414000000000003
$ws.Range("C10").Value = 0.012661
$ws.Range("C11").Value = 4.1265000000000001
$ws.Range("C12").Value = 5.2393000000000001
$ws.Range("C13").Value = 0.13689999999999999
$ws.Range("C14").Value = 0.034842999999999999
$ws.Range("C15").Value = 0.1749
$ws.Range("C16").Value = 0.61060000000000003
$ws.Range("C17").Value = 0.030074999999999998
$ws.Range("C18").Value = 0.44879999999999998
$ws.Range("C19").Value = 0.4466
$ws.Range("C20").Value = 0.59960000000000002
$ws.Range("C21").Value = 0.92259999999999998
$ws.Range("C22").Value = 2.3220000000000001
$ws.Range("C23").Value = 0.4708
$ws.Range("C24").Value = 1.1598999999999999
$ws.Range("C25").Value = 0.0053949999999999996
$ws.Range("C26").Value = 0.078700000000000006
$ws.Range("C27").Value = 0.18940000000000001
$ws.Range("C28").Value = 0.26190000000000002
$ws.Range("C29").Value = 0.68620000000000003
$ws.Range("C30").Value = 0.92400000000000004
$ws.Range("C31").Value = 0.049500000000000002
$ws.Range("C32").Value = 0.00026233000000000002
$ws.Range("C33").Value = 0.051060999999999995
$ws.Range("C34").Value = 0.0034129999999999998
$ws.Range("C35").Value = 0.58250000000000002
$ws.Range("C36").Value = 5.4436

$ws.Range("A4").Select()
